$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2844.6428
$ws.Range("I43").Value = 2740
$ws.Range("J43").Value = 2886.5
$ws.Range("K43").Value = 2740
$ws.Range("L43").Value = 2886.5
$ws.Range("M43").Value = -2671
$ws.Range("N43").Value = -3024.5

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H98").Value = 17809082
$ws.Range("I98").Value = 21045854
$ws.Range("J98").Value = 6832.5
$ws.Range("K98").Value = 21045854
$ws.Range("L98").Value = 6832.5
$ws.Range("M98").Value = -21044356
$ws.Range("N98").Value = -9828.5

$ws.Range("H121").Value = 1141.5518
$ws.Range("J121").Value = 1128.75
$ws.Range("L121").Value = 3386.25
$ws.Range("N121").Value = -6880.25

$ws.Range("H122").Value = 17809082
$ws.Range("I122").Value = 21045854
$ws.Range("J122").Value = 6832.5
$ws.Range("K122").Value = 63137562
$ws.Range("L122").Value = 20497.5
$ws.Range("M122").Value = -63135112
$ws.Range("N122").Value = -25397.5

$ws.Range("H123").Value = 32900
$ws.Range("J123").Value = 32900
$ws.Range("L123").Value = 32900
$ws.Range("N123").Value = -42700

$ws.Range("H124").Value = 39001
$ws.Range("J124").Value = 39001
$ws.Range("L124").Value = 39001
$ws.Range("N124").Value = -48821

$ws.Range("H126").Value = 30217.705
$ws.Range("J126").Value = 30217.705
$ws.Range("L126").Value = 30217.705
$ws.Range("N126").Value = -40097.705

$ws.Range("H128").Value = 49982.316
$ws.Range("J128").Value = 49982.316
$ws.Range("L128").Value = 49982.316
$ws.Range("N128").Value = -59942.316

$ws.Range("H130").Value = 50000
$ws.Range("J130").Value = 50000
$ws.Range("L130").Value = 50000
$ws.Range("N130").Value = -60040

$ws.Range("H137").Value = 190174.25
$ws.Range("I137").Value = 341237.34
$ws.Range("J137").Value = 1345.3928
$ws.Range("K137").Value = 1023712.02
$ws.Range("L137").Value = 4036.1784
$ws.Range("M137").Value = -1021162.02
$ws.Range("N137").Value = -9136.178400000001

$ws.Range("H138").Value = 3509.8374
$ws.Range("I138").Value = 5546.0605
$ws.Range("J138").Value = 2080.149
$ws.Range("K138").Value = 16638.1815
$ws.Range("L138").Value = 6240.447
$ws.Range("M138").Value = -11498.1815
$ws.Range("N138").Value = -16520.447

$ws.Range("H141").Value = 9579.286
$ws.Range("I141").Value = 11456.667
$ws.Range("K141").Value = 34370.001
$ws.Range("M141").Value = -29190.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13036.573
$ws.Range("I32").Value = 9539.712
$ws.Range("J32").Value = 28991
$ws.Range("K32").Value = 9539.712
$ws.Range("L32").Value = 28991
$ws.Range("M32").Value = -9252.712
$ws.Range("N32").Value = -29565

$ws.Range("H61").Value = 2925.8572
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2925.8572
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2925.8572
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -3349.8572

$ws.Range("H74").Value = 1366.1666
$ws.Range("I74").Value = 766.6667
$ws.Range("J74").Value = 1965.6666
$ws.Range("K74").Value = 766.6667
$ws.Range("L74").Value = 1965.6666
$ws.Range("M74").Value = 107.3333
$ws.Range("N74").Value = -3713.6666

$ws.Range("H77").Value = 1366.1666
$ws.Range("I77").Value = 766.6667
$ws.Range("J77").Value = 1965.6666
$ws.Range("K77").Value = 3833.3335
$ws.Range("L77").Value = 9828.333000000001
$ws.Range("M77").Value = 534.6665000000003
$ws.Range("N77").Value = -18564.333

$ws.Range("H132").Value = 3474349.2
$ws.Range("I132").Value = 5953665
$ws.Range("J132").Value = 3307.0667
$ws.Range("K132").Value = 17860995
$ws.Range("L132").Value = 9921.2001
$ws.Range("M132").Value = -17858465
$ws.Range("N132").Value = -14981.2001

$ws.Range("H136").Value = 2925.8572
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 2925.8572
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 8777.571599999999
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -13877.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3109.4546
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 3109.4546
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 9328.363799999999
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -14398.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2250778.2
$ws.Range("I58").Value = 4360583.5
$ws.Range("J58").Value = 4856.516
$ws.Range("K58").Value = 4360583.5
$ws.Range("L58").Value = 4856.516
$ws.Range("M58").Value = -4360380.5
$ws.Range("N58").Value = -5262.516

$ws.Range("H132").Value = 7411506.5
$ws.Range("I132").Value = 13334561
$ws.Range("J132").Value = 7688.2
$ws.Range("K132").Value = 40003683
$ws.Range("L132").Value = 23064.6
$ws.Range("M132").Value = -40001153
$ws.Range("N132").Value = -28124.6

$ws.Range("H134").Value = 3473970
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 3473970
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 10421910
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -10426980

$ws.Range("H136").Value = 2250778.2
$ws.Range("I136").Value = 4360583.5
$ws.Range("J136").Value = 4856.516
$ws.Range("K136").Value = 13081750.5
$ws.Range("L136").Value = 14569.548
$ws.Range("M136").Value = -13079200.5
$ws.Range("N136").Value = -19669.548

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 769.3022999999999
$ws.Range("I5").Value = 280.35
$ws.Range("J5").Value = 1194.4783
$ws.Range("K5").Value = 841.0500000000001
$ws.Range("L5").Value = 3583.4349
$ws.Range("M5").Value = -729.0500000000001
$ws.Range("N5").Value = -3807.4349

$ws.Range("H75").Value = 1924.375
$ws.Range("I75").Value = 1490
$ws.Range("J75").Value = 2069.1667
$ws.Range("K75").Value = 4470
$ws.Range("L75").Value = 6207.500100000001
$ws.Range("M75").Value = -3472
$ws.Range("N75").Value = -8203.500100000001

$ws.Range("H78").Value = 1924.375
$ws.Range("I78").Value = 1490
$ws.Range("J78").Value = 2069.1667
$ws.Range("K78").Value = 13410
$ws.Range("L78").Value = 18622.5003
$ws.Range("M78").Value = -8418
$ws.Range("N78").Value = -28606.5003

$ws.Range("H135").Value = 769.3022999999999
$ws.Range("I135").Value = 280.35
$ws.Range("J135").Value = 1194.4783
$ws.Range("K135").Value = 2523.15
$ws.Range("L135").Value = 10750.3047
$ws.Range("M135").Value = 11.84999999999991
$ws.Range("N135").Value = -15820.3047

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H43").Value = 881
$ws.Range("I43").Value = 881
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 881
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -730
$ws.Range("N43").ClearContents()

$ws.Range("H46").Value = 25000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 25000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 25000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -25312

$ws.Range("H57").Value = 5506.75
$ws.Range("J57").Value = 5999.857
$ws.Range("L57").Value = 5999.857
$ws.Range("N57").Value = -7639.857

$ws.Range("H80").Value = 2463.9
$ws.Range("I80").Value = 2076.5
$ws.Range("J80").Value = 3045
$ws.Range("K80").Value = 2076.5
$ws.Range("L80").Value = 3045
$ws.Range("M80").Value = -1078.5
$ws.Range("N80").Value = -5041

$ws.Range("H83").Value = 2463.9
$ws.Range("I83").Value = 2076.5
$ws.Range("J83").Value = 3045
$ws.Range("K83").Value = 10382.5
$ws.Range("L83").Value = 15225
$ws.Range("M83").Value = -5390.5
$ws.Range("N83").Value = -25209

$ws.Range("H122").Value = 142859520
$ws.Range("I122").Value = 200002850
$ws.Range("J122").Value = 1195
$ws.Range("K122").Value = 600008550
$ws.Range("L122").Value = 3585
$ws.Range("M122").Value = -600006100
$ws.Range("N122").Value = -8485

$ws.Range("H126").Value = 1957
$ws.Range("I126").Value = 1301.3334
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 3904.0002
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -1434.0002
$ws.Range("N126").Value = -13340

$ws.Range("H132").Value = 21301240
$ws.Range("I132").Value = 33368218
$ws.Range("J132").Value = 6576.294
$ws.Range("K132").Value = 100104654
$ws.Range("L132").Value = 19728.882
$ws.Range("M132").Value = -100102124
$ws.Range("N132").Value = -24788.882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3372.1
$ws.Range("I40").Value = 3546.8333
$ws.Range("K40").Value = 3546.8333
$ws.Range("M40").Value = -3410.8333

$ws.Range("H122").Value = 113647370
$ws.Range("I122").Value = 166684910
$ws.Range("J122").Value = 50002300
$ws.Range("K122").Value = 500054730
$ws.Range("L122").Value = 150006900
$ws.Range("M122").Value = -500052280
$ws.Range("N122").Value = -150011800

$ws.Range("H132").Value = 3078840
$ws.Range("I132").Value = 7144881.5
$ws.Range("J132").Value = 1835.3784
$ws.Range("K132").Value = 21434644.5
$ws.Range("L132").Value = 5506.135200000001
$ws.Range("M132").Value = -21432114.5
$ws.Range("N132").Value = -10566.1352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14527667
$ws.Range("I136").Value = 7142334
$ws.Range("K136").Value = 21427002
$ws.Range("L136").Value = 166671888
$ws.Range("M136").Value = -21424452
$ws.Range("N136").Value = -166676988

